$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark Section 8 (Programming) videos as completed (G column = TARGET/COMPLETED date)
$ws.Range("G86:G107").Value = 42794

# Shift TARGET DATE for Section 9 (Materials and Effects) from Feb 28 to Mar 1, 2017
$ws.Range("H109:H125").Value = 42795

# Shift TARGET DATE for Section 12 (Audio) rows from Mar 1 to Mar 2, 2017
$ws.Range("H165:H172").Value = 42796
$ws.Range("H174:H182").Value = 42796

# Shift TARGET DATE for Section 18/19 rows from Mar 2 to Mar 3, 2017
$ws.Range("H225").Value = 42797
$ws.Range("H227:H243").Value = 42797

# Reduce Mock Exam quiz durations from 5 to 3 minutes
$ws.Range("C227:C242").Value = 3

# Update selection to match author state (cursor moved to G107 after marking
# section 8 complete). Note: the frozen-pane scroll anchor (topLeftCell) is
# view-state that this COM host does not expose a settable property for
# (no Window/Pane object in the bridged model), so only the selection can be
# reproduced here.
$ws.Range("G107").Select()
